# data : case 1
# Swap the two column widths (col A <-> col B) and update the numeric
# values in A1:B3 (row 4 is left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap column widths -------------------------------------------------
# Original: col A width = 14.7109375, col B width = 15.42578125
# Target:   col A width = 15.42578125, col B width = 14.7109375
$ws.Columns.Item(1).ColumnWidth = 14.66
$ws.Columns.Item(2).ColumnWidth = 13.8

# --- Update cell values --------------------------------------------------
$ws.Range("A1").Value = 0.050267379504308984
$ws.Range("B1").Value = -0.050267379563267697

$ws.Range("A2").Value = -0.0077125121154475052
$ws.Range("B2").Value = 0.0077125120165095996

$ws.Range("A3").Value = -0.051993854651702823
$ws.Range("B3").Value = 0.051993854582379387
